$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Robot designs" becomes "Part" ---
$ws.Range("A1").Value = "Part"

# --- Column widths (B and D get wider) ---
# ColumnWidth setter rounds to the nearest pixel internally (raw XML width =
# input + 5/6, then snapped to 1/6ths), so feed it target-minus-padding to
# land on the closest achievable value.
$ws.Columns(2).ColumnWidth = 37.830729166666664
$ws.Columns(4).ColumnWidth = 48.608072916666664

# --- New row 2 content ---
$ws.Range("A2").Value = "Actuators"

$desc = "1) High torque density" + [char]10 + "2) Force transparency" + [char]10 + "3) Mechanical robustness" + [char]10 + "4) Energy efficieny" + [char]10 + "5) Low mechanical impendance"
$ws.Range("B2").Value = $desc
$ws.Range("B2").WrapText = $true

# Hyperlink to the thesis PDF. TextToDisplay seeds the hyperlink element's
# "display" attribute (kept as the relative path); the cell's visible text is
# then overwritten to just the filename, matching the saved workbook.
$ws.Hyperlinks.Add($ws.Range("D2"), "../Bibliography/TotalBibliography/PetrosPolidorouDiplomaThesisMotorDesign.pdf", "", "", "../Bibliography/TotalBibliography/PetrosPolidorouDiplomaThesisMotorDesign.pdf")
$ws.Range("D2").Value = "PetrosPolidorouDiplomaThesisMotorDesign.pdf"

# Row height grows to fit the wrapped, multi-line description.
$ws.Rows(2).RowHeight = 130.19999999999999

# --- Selection moves to B9 ---
$null = $ws.Range("B9").Select()
